$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.173.22'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.886.38'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'321.55"
$ws.Range("E5").Value = '  -3.13%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'0.4704"
$ws.Range("E7").Value = '  +2.12%  '
$ws.Range("D8").Value = "'0.4020"
$ws.Range("E8").Value = '  -2.80%  '
$ws.Range("D9").Value = "'47.26"
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = "'0.9928"
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D12").Value = "'22.38"
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '1.872.00'
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").Value = "'5.883"
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").Value = "'7.029"
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = "'88.75"
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").Value = "'0.00001019"
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = "'17.43"
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = '29.181.57'
$ws.Range("D23").Value = "'5.481"
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").Value = "'11.69"
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("D25").Value = "'2.176"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("D26").Value = '2.116.61'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Value = "'154.85"
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  +6.23%  '
$ws.Range("D30").Value = "'2.073"
$ws.Range("E30").Value = '  -2.46%  '
$ws.Range("D31").Value = "'117.24"
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").Value = "'0.09435"
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("D34").Value = "'3.536"
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").Value = "'1.379"
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("D36").Value = "'5.348"
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = "'0.06065"
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").Value = "'0.02225"
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").Value = "'1.173"
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").Value = "'8.014"
$ws.Range("E40").Value = '  -5.13%  '
$ws.Range("D41").Value = "'0.5810"
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = "'0.1825"
$ws.Range("D43").Value = "'2.479"
$ws.Range("E43").Value = '  +7.30%  '
$ws.Range("D44").Value = "'10.01"
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").Value = "'1.274"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").Value = "'0.07688"
$ws.Range("E46").Value = '  +2.12%  '
$ws.Range("D47").Value = "'12.08"
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("D48").Value = "'0.5459"
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").Value = "'1.899"
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("D50").Value = "'113.36"
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").Value = "'44.30"
$ws.Range("E51").Value = '  +0.13%  '
